$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# PHASE 1: insert 5 new blank "Title and Content" slides right after slide 5
# (these become the new section-divider slides for ids 266,265,267,268,269)
# ---------------------------------------------------------------------------
$null = $p.Slides.Add(6, 2)
$null = $p.Slides.Add(7, 2)
$null = $p.Slides.Add(8, 2)
$null = $p.Slides.Add(9, 2)
$null = $p.Slides.Add(10, 2)

# ---------------------------------------------------------------------------
# PHASE 2: titles for the 5 new slides (bodies are already empty placeholders)
# ---------------------------------------------------------------------------
$p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange.Text = "Starting sagas"
$p.Slides.Item(7).Shapes.Item(1).TextFrame.TextRange.Text = "Storing state"
$p.Slides.Item(8).Shapes.Item(1).TextFrame.TextRange.Text = "Mapping messages to sagas"
$p.Slides.Item(9).Shapes.Item(1).TextFrame.TextRange.Text = "Requesting timeouts"
$p.Slides.Item(10).Shapes.Item(1).TextFrame.TextRange.Text = "Emitting messages"

# ---------------------------------------------------------------------------
# PHASE 3: slide 5 -- repurposed from "Sagas = message driven state machines"
# to "Handling messages" with an empty body. Recreate it so the leftover
# "* " bullet paragraph formatting doesn't linger in the (now empty) body.
# ---------------------------------------------------------------------------
$p.Slides.Item(5).Delete()
$s5 = $p.Slides.Add(5, 2)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Handling messages"

# ---------------------------------------------------------------------------
# PHASE 4: slide 4 ("Saga definition") -- append two paragraphs to the body
# ---------------------------------------------------------------------------
$s4Body = $p.Slides.Item(4).Shapes.Item(2).TextFrame.TextRange
$s4Body.Text = $s4Body.Text + "`r`rSagas = message driven state machines"

# ---------------------------------------------------------------------------
# PHASE 5: the original slides 6,7,8,9 are now at positions 11,12,13,14
# ---------------------------------------------------------------------------

# 11: "Sample domain walkthrough" -- unchanged, no edits needed

# 12: "Exercise 1" -> "Exercise 1 - Order Policy saga", body rebuilt
$s12 = $p.Slides.Item(12)
$s12.Shapes.Item(1).TextFrame.TextRange.Text = "Exercise 1 " + [char]0x2013 + " Order Policy saga"
$s12Body = $s12.Shapes.Item(2).TextFrame.TextRange
$s12Body.Text = ""
$s12Body.Text = "State changes`rStartOrder (Command)`rPlaceOrder (Command)`rCancelOrder (Command)`rOrderAbandoned (Event)`rBusiness rules:`rAn order is abandoned if not cancelled or placed within 20 seconds`rEvents should be emitted for each relevant state change"
$s12Body.Paragraphs(2,1).IndentLevel = 2
$s12Body.Paragraphs(3,1).IndentLevel = 2
$s12Body.Paragraphs(4,1).IndentLevel = 2
$s12Body.Paragraphs(5,1).IndentLevel = 2
$s12Body.Paragraphs(7,1).IndentLevel = 2
$s12Body.Paragraphs(8,1).IndentLevel = 2

# 13: "Walkthrough - Exercise 1" -> "Walkthrough" (body "Exercise 1" unchanged)
$s13 = $p.Slides.Item(13)
$s13.Shapes.Item(1).TextFrame.TextRange.Text = "Walkthrough"

# 14: "Upcoming saga changes in v6" -- unchanged, no edits needed
